$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version (row 3), Status (row 6), Date (row 8) and Publisher's Contact (row 10)
$ws.Cells.Item(3, 2).Value = "0.1.7"
$ws.Cells.Item(6, 2).Value = "draft"
$ws.Cells.Item(8, 2).Value = "2024-11-22T12:33:30-06:00"
$ws.Cells.Item(10, 2).Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# Insert a new row 11 for the author's Contact, copying the formatting of the
# existing Contact row (row 10) so border/fill/font stay identical. This
# pushes the old "Description" row (and everything after it) down by one;
# the row that lands on 12 is a leftover duplicate of the old Contact row
# (rows 10 and 11 used to share the same text), which we then repurpose
# below as the new "Jurisdiction" row -- so no second insert is needed.
$ws.Rows.Item(11).Insert()
$ws.Range("A10:B10").Copy($ws.Range("A11:B11"))
$ws.Cells.Item(11, 1).Value = "Contact"
$ws.Cells.Item(11, 2).Value = "Bob Milius (bmilius@nmdp.org)"

$ws.Cells.Item(12, 1).Value = "Jurisdiction"
$ws.Cells.Item(12, 2).Value = ""
